# ---------------------------------------------------------------------------
# Applies the "aula11 / SolveProblem" edit described in the commit:
#   - rename sheet "Plan1" -> "SolveProblem"
#   - rebuild its small "problem solving steps" table with new content
#   - change its zoom from 175% to 115%
#   - make "Folha2021" the active sheet, with A1:S1 (the merged title) selected
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename "Plan1" to "SolveProblem"
# ---------------------------------------------------------------------------
$solve = $wb.Worksheets.Item("Plan1")
$solve.Name = "SolveProblem"

# ---------------------------------------------------------------------------
# 2) Rebuild the sheet's content: clear the old "Problema" helper table
#    (which lived in D6:G12) and write the new "Etapas" table in A1:D7.
# ---------------------------------------------------------------------------
$solve.Cells.Clear() | Out-Null

# Header row (bold)
$solve.Range("A1").Value = "Etapas"
$solve.Range("B1").Value = "Solução de"
$solve.Range("C1").Value = "Problema"
$solve.Range("D1").Value = "O que é problema"

# Row 2
$solve.Range("A2").Value = 0
$solve.Range("B2").Value = "Entender o Problema"
$solve.Range("C2").Value = "Tem solução Trivial"
$solve.Range("D2").Value = "Não é problema"

# Row 3
$solve.Range("C3").Value = "Tem solução não Trivial"
$solve.Range("D3").Value = "É problema"

# Row 4
$solve.Range("C4").Value = "Não tem solução"
$solve.Range("D4").Value = "Não é problema"

# Row 5
$solve.Range("A5").Value = 1
$solve.Range("B5").Value = "Decompor"
$solve.Range("C5").Value = "Seprara em problemas menores"
$solve.Range("D5").Value = "Triviais"

# Row 6
$solve.Range("A6").Value = 2
$solve.Range("B6").Value = "Encontrar padrões"

# Row 7
$solve.Range("A7").Value = 3
$solve.Range("B7").Value = "Criar um algoritmo"
$solve.Range("C7").Value = "Receita, de passos"

# Bold header row, like the rest of the workbook's table headers
$solve.Range("A1:D1").Font.Bold = $true

# Column widths roughly matching the auto-fit content widths of the new table
$solve.Columns.Item(1).ColumnWidth = 5.877604166666667
$solve.Columns.Item(2).ColumnWidth = 19.166666666666668
$solve.Columns.Item(3).ColumnWidth = 29.022135416666668
$solve.Columns.Item(4).ColumnWidth = 16.166666666666668

# Page setup (A4 / portrait) for the sheet
$solve.PageSetup.PaperSize = 9
$solve.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 3) Zoom out from 175% to 115% while this sheet is active
# ---------------------------------------------------------------------------
$solve.Activate()
$excel.ActiveWindow.Zoom = 115

# ---------------------------------------------------------------------------
# 4) Switch to "Folha2021" and select the merged title range A1:S1
# ---------------------------------------------------------------------------
$folha = $wb.Worksheets.Item("Folha2021")
$folha.Activate()
$folha.Range("A1:S1").Select()
